$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDF_CH_sub_1_stratification")

# Row 2: race changes from White -> Native Hawaiian or Other Pacific Islander
$ws.Range("C2").Value = "Native Hawaiian or Other Pacific Islander"

# Row 3: ethnicity changes Not Hispanic or Latino -> Hispanic or Latino
#        race changes White -> Native Hawaiian or Other Pacific Islander
$ws.Range("B3").Value = "Hispanic or Latino"
$ws.Range("C3").Value = "Native Hawaiian or Other Pacific Islander"

# Row 5: ethnicity changes Hispanic or Latino -> Not Hispanic or Latino
#        race changes Other Race -> Native Hawaiian or Other Pacific Islander
$ws.Range("B5").Value = "Not Hispanic or Latino"
$ws.Range("C5").Value = "Native Hawaiian or Other Pacific Islander"
